# Cover page update: bump the release version and date shown under the
# document title.
#
#   "Version 11.08.01, 2016-02-15"  ->  "Version 11.10.01, 2016-05-02"
#
# (matches the commit message: "preparing for 11.10.01 release")
#
# Four single characters change in the underlying OOXML runs:
#   "0" -> "1"   (tens digit of the minor version: 08 -> 10)
#   "8" -> "0"
#   "2" -> "5"   (month of the date: 02 -> 05)
#   "15" -> "02" (day of the date: 15 -> 02)
#
# Replace the date first, then the version, so each Find/Replace only
# touches the text it needs to and the surrounding "Version 11." runs
# are left untouched.

$d = $word.ActiveDocument

$d.Content.Find.Execute("2016-02-15", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2016-05-02", 2)

$d.Content.Find.Execute("08.01", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "10.01", 2)
